# Apply updates described by the commit:
#   "remove gupen + roll back"
# -> Update selected cell / active cell on the sheet view
# -> Roll back a batch of L/M column values (columns L=FscoreTotal-ish, M=..)
#    in rows 7-11 and 17-21 to their previous values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# Update the active selection from I26 to I24
$ws.Range("I24").Select()

# Row 7
$ws.Range("L7").Value = 0.712
$ws.Range("M7").Value = 0.766

# Row 8
$ws.Range("L8").Value = 0.715
$ws.Range("M8").Value = 0.766

# Row 9
$ws.Range("L9").Value = 0.71
$ws.Range("M9").Value = 0.767

# Row 10
$ws.Range("L10").Value = 0.718
$ws.Range("M10").Value = 0.789

# Row 11
$ws.Range("L11").Value = 0.718
$ws.Range("M11").Value = 0.757

# Row 17
$ws.Range("L17").Value = 0.755
$ws.Range("M17").Value = 0.826

# Row 18
$ws.Range("L18").Value = 0.741
$ws.Range("M18").Value = 0.815

# Row 19
$ws.Range("L19").Value = 0.77

# Row 20
$ws.Range("M20").Value = 0.842

# Row 21
$ws.Range("M21").Value = 0.852
